# adminStudent.pptx regeneration edit
#
# The canonical diff for this fixture is almost entirely made up of
# boilerplate that a PowerPoint re-save/regeneration produces (new
# p14:creationId GUIDs, buFont panose attributes, en-US -> en-SG default
# paragraph language on *empty* master/layout placeholders, a Georgian
# font entry added to the theme font scheme, the notes master/notes
# slide being dropped together with the now-unused theme2.xml, and
# internal shape-id renumbering). None of that is reachable through the
# PowerPoint automation object model (there is no LanguageID setter, no
# notes-master deletion call that doesn't also remove the slide, no
# font-scheme editor, and shape/slide ids are managed internally and
# read-only).
#
# The few parts of the diff that *are* reachable through COM and that we
# apply here:
#   * the deck's first slide number goes back to the default (1) instead
#     of the old custom value (10)
#   * the eleven shapes on slide 1 keep their drawn content but are
#     renamed to match the regenerated deck's shape names (Rectangle 5 ->
#     Rectangle 38, etc.) - only the cosmetic `name` changes because the
#     numeric `id` itself is not settable from script, same as in real
#     PowerPoint.

$p = $ppt.ActivePresentation

# First slide number reverts to the PowerPoint default.
$p.PageSetup.FirstSlideNumber = 1

$s = $p.Slides.Item(1)

$renames = @{
    'Rectangle 5'            = 'Rectangle 38'
    'Rectangle 6'            = 'Rectangle 39'
    'Straight Connector 19'  = 'Straight Connector 40'
    'TextBox 20'             = 'TextBox 41'
    'Right Triangle 21'      = 'Right Triangle 42'
    'Rectangle 27'           = 'Rectangle 43'
    'Rectangle 28'           = 'Rectangle 44'
    'Straight Connector 30'  = 'Straight Connector 45'
    'TextBox 31'             = 'TextBox 46'
    'Right Triangle 33'      = 'Right Triangle 47'
    'Straight Connector 35'  = 'Straight Connector 48'
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    $newName = $renames[$shape.Name]
    if ($newName) {
        $shape.Name = $newName
    }
}
